$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# The sheet originally has 2 data rows (row 2 and row 3). This edit:
#   1) swaps the contents of row 2 and row 3
#   2) appends a new row 4 with a new observation
# ------------------------------------------------------------------

# Preserve the "empty" marker cells (I, K, AT, AY) for the new row by
# copying the already-empty cells from row 2 - this keeps those
# columns present (but blank) on row 4, matching the other rows.
$ws.Range("I2").Copy($ws.Range("I4"))
$ws.Range("K2").Copy($ws.Range("K4"))
$ws.Range("AT2").Copy($ws.Range("AT4"))
$ws.Range("AY2").Copy($ws.Range("AY4"))

# --- Row 2 becomes what row 3 used to hold ---
$ws.Range("A2").Value = 112331171
$ws.Range("B2").Value = 90800
$ws.Range("D2").Value = "NT"
$ws.Range("E2").Value = 3100
$ws.Range("F2").Value = "Talltaggsvamp"
$ws.Range("G2").Value = "Bankera fuligineoalba"
$ws.Range("H2").Value = "(Schmidt : Fr.) Pouzar"
$ws.Range("R2").Value = 7537240

# --- Row 3 becomes what row 2 used to hold ---
$ws.Range("A3").Value = 112331579
$ws.Range("B3").Value = 88203
$ws.Range("D3").Value = "VU"
$ws.Range("E3").Value = 6286
$ws.Range("F3").Value = "Torrmusseron"
$ws.Range("G3").Value = "Tricholoma sudum"
$ws.Range("H3").Value = "(Fr.) Quél."
$ws.Range("R3").Value = 7537233

# --- New row 4 ---
$ws.Range("A4").Value = 112534860
$ws.Range("B4").Value = 90808
$ws.Range("C4").Value = "Ovaliderad"
$ws.Range("D4").Value = "NT"
$ws.Range("E4").Value = 4362
$ws.Range("F4").Value = "Blå taggsvamp"
$ws.Range("G4").Value = "Hydnellum caeruleum"
$ws.Range("H4").Value = "(Hornem.) P.Karst."
$ws.Range("P4").Value = "Esrangeåsen (Esrangeåsen), T lm"
$ws.Range("Q4").Value = 748912
$ws.Range("R4").Value = 7537606
$ws.Range("S4").Value = 10
$ws.Range("T4").Value = "Norrbotten"
$ws.Range("U4").Value = "Kiruna"
$ws.Range("V4").Value = "Torne lappmark"
$ws.Range("W4").Value = "Jukkasjärvi"

# Dates are stored as plain text (as in rows 2/3), not as Excel date
# serials, so force text formatting before assigning, then strip the
# format again so no stray style is left behind on the cell.
$ws.Range("Y4").NumberFormat = "@"
$ws.Range("Y4").Value = "2023-10-05"
$ws.Range("Y4").ClearFormats()

$ws.Range("AA4").NumberFormat = "@"
$ws.Range("AA4").Value = "2023-10-05"
$ws.Range("AA4").ClearFormats()

$ws.Range("AD4").Value = $false
$ws.Range("AE4").Value = $false
$ws.Range("AG4").Value = $false
$ws.Range("AW4").Value = "per-erik mukka"
$ws.Range("AX4").Value = "per-erik mukka"
